# Week 8 journal rewrite: re-centers/re-styles the title & "Simplicity" heading,
# drops the old tab-driven "normal0" paragraph style in favor of direct
# character formatting (Arial east-asian font, cstheme minorHAnsi, black color),
# and edits/merges the body paragraphs' wording per the journal update.
#
# The formatting/run-splitting changes are intricate (styles removed, runs
# merged/split, paragraph marks re-styled), so the most reliable way to apply
# them exactly is to replace the whole document body with the target WordML
# via Range.InsertXML - this is a normal Word COM operation (it replaces only
# the content of the Range it's called on; here that's $d.Content, i.e. the
# body up to - but not including - the final section-properties mark, so the
# sectPr/page setup at the end of the document is left untouched).

$d = $word.ActiveDocument

$bodyXml = @'
<w:p>
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Arial" w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:color w:val="000000"/>
          <w:sz w:val="32"/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Arial" w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:color w:val="000000"/>
          <w:sz w:val="32"/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>T</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Arial" w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:color w:val="000000"/>
          <w:sz w:val="32"/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>eam Project  - Journal</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:rFonts w:eastAsia="Arial" w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:color w:val="000000"/>
          <w:sz w:val="32"/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Arial" w:cstheme="minorHAnsi"/>
          <w:b/>
          <w:color w:val="000000"/>
          <w:sz w:val="32"/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Simplicity</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:eastAsia="Arial" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:eastAsia="Arial" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Arial" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t xml:space="preserve">In this week’s meeting our team decided to add design patterns in our project to simplify the functionalities and to make sure loose coupling of logic with data. Each team member researched about the patterns and where they could be applied in efficient way. Finally, each team member picked different pattern which when applied server it’s purpose of increasing simplicity and </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Arial" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>making code easy to understand.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:eastAsia="Arial" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Arial" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>I worked on and implemented State Machine pattern. This pattern is applied to figure out status of ship(player) at any point of time in the game. The ship could be on Ordinary Island, Treasure Island or travelling. Accordingly, logic is implemented which is appropriate to the state the ship is in. This simplifies the code in a way that, class is responsible and only need to implement only event occurring when ship in that state.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:eastAsia="Arial" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Arial" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>We are also making sure that each member applies pattern in such a way that, the resulting code become compatible with rest of the code of the game. We believe that after applying these patterns our code will become more robust, coherent and simple. Additionally, we have also worked on implementing timer logic to make a timer based game in which multiple players can join and play game for certain time and at the end of the time, winner(s) w</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Arial" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>ill be displayed by the server.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:rFonts w:cstheme="minorHAnsi"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Arial" w:cstheme="minorHAnsi"/>
          <w:color w:val="000000"/>
        </w:rPr>
        <w:t>We will be  testing our code now to ensure stable running in all scenarios along with minor changes.</w:t>
      </w:r>
    </w:p>
'@

$xmlDoc = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>$bodyXml</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$d.Content.InsertXML($xmlDoc)

Write-Host "Updated document. Paragraph count:" $d.Paragraphs.Count
